# Updated cryptos list on Mon Aug 14 11:44:16 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (NumberFormat "@") before assignment so purely-numeric-looking
# strings (e.g. "0.9990", "0.00001055") are preserved verbatim instead of being
# auto-converted to numeric values by Excel, matching the source inline-string cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.338.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.01%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07482'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.45%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.48'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07738'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.844.22'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.984'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6796'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001055'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.94'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.174'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.361.15'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.85'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.19%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9992'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.503'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.36%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9990'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.59'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.427'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.06%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.97%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06580'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +17.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.416'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.480'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.104'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.089'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.15%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6952'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.581'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.263.18'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.03%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.34%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.794'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9159'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9990'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.007.80'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.52%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.11'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.731'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.072'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1160'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.953'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3947'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.93%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.06%  '
